$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.859.25'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '3.460.05'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").Value = '180.33'
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("E7").Value = '  +2.39%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '3.459.95'
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").Value = '6.98'
$ws.Range("E11").Value = '  -2.28%  '
$ws.Range("D12").Value = '0.428'
$ws.Range("E12").Value = '  -2.26%  '
$ws.Range("D13").Value = '4.061.27'
$ws.Range("E13").Value = '  -1.43%  '
$ws.Range("D14").Value = '32.08'
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("D16").Value = '67.848.92'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("E17").Value = '  -3.46%  '
$ws.Range("D18").Value = '3.462.93'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").Value = '6.19'
$ws.Range("D20").Value = '14.08'
$ws.Range("E20").Value = '  -5.24%  '
$ws.Range("D21").Value = '391.96'
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("D22").Value = '7.89'
$ws.Range("E22").Value = '  -2.39%  '
$ws.Range("E23").Value = '  +2.36%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").Value = '71.81'
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("E27").Value = '  -4.94%  '
$ws.Range("D28").Value = '10.40'
$ws.Range("E28").Value = '  -2.81%  '
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = '6.10'
$ws.Range("E31").Value = '  -3.25%  '
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("E33").Value = '  -5.80%  '
$ws.Range("D34").Value = '23.46'
$ws.Range("E34").Value = '  -3.07%  '
$ws.Range("E35").Value = '  -2.08%  '
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("E37").Value = '  -7.36%  '
$ws.Range("D38").Value = '162.02'
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("D39").Value = '0.887'
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("E41").Value = '  -5.23%  '
$ws.Range("D42").Value = '4.64'
$ws.Range("E42").Value = '  -2.80%  '
$ws.Range("D43").Value = '6.69'
$ws.Range("E43").Value = '  -7.29%  '
$ws.Range("D44").Value = '26.08'
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("D45").Value = '0.0719'
$ws.Range("E45").Value = '  -3.29%  '
$ws.Range("D46").Value = '26.12'
$ws.Range("E46").Value = '  -6.65%  '
$ws.Range("D47").Value = '2.724.99'
$ws.Range("E47").Value = '  -4.76%  '
$ws.Range("D48").Value = '41.25'
$ws.Range("E48").Value = '  -2.61%  '
$ws.Range("D49").Value = '0.0297'
$ws.Range("E49").Value = '  -3.18%  '
$ws.Range("D50").Value = '328.56'
$ws.Range("E50").Value = '  -5.93%  '
$ws.Range("E51").Value = '  -5.24%  '
